# Adicionando esquema 4-3-3 para jogadores com maiores scores e atualizando tabela resumo
#
# The sheet already has a "SCORE 5-3-2 - Time de Maior Score (media)" summary
# row at row 16. We insert a brand new row above it for the "SCORE 4-3-3 -
# Time de Maior Score (media)" summary (which pushes the existing 5-3-2 row
# down to row 17, formulas/relative references adjust automatically), fill
# in the new row's data, and fix up the row that got pushed down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 16 (old row 16 becomes row 17) -------
$ws.Rows.Item(16).Insert()

# --- 2. Populate the new row 16 : "SCORE 4-3-3 - Time de Maior Score (media)"
$ws.Range("B16").Value = "SCORE 4-3-3 - Time de Maior Score (media)"

$newRow16Values = @(133.57999999999998,148.14000000000001,110.25,130.72,117.98,136.13000000000002,95.14,115.84000000000002,101.12000000000002,132.85000000000002,85.26,74.81,78.64,85.99,73.8,90.420000000000016,69.77000000000001,99.39,76.330000000000013,63.39,35.160000000000004,51.999999999999993,50.25,44.589999999999996,70.290000000000006,56.59,37.92,46.84,26.740000000000002,55.2,74.62,37.950000000000003,37.359999999999992,62.160000000000004,54.66,57.550000000000004,110.76000000000003,88.350000000000009)

for ($i = 0; $i -lt $newRow16Values.Length; $i++) {
    $ws.Cells.Item(16, 3 + $i).Value = $newRow16Values[$i]
}

$ws.Range("AO16").Formula = "=SUM(C16:AN16)"
$ws.Range("AQ16").Formula = "=(AO16*100)/`$AO`$6"

# --- 3. Row 17 (formerly row 16, "SCORE 5-3-2 - Time de Maior Score (media)")
#        already kept its values/formulas from the insert/shift. It just
#        needs the extra AP17 placeholder cell (formatted like its
#        neighbours) to match the reference layout.
$ws.Range("AN17").Copy()
$ws.Range("AP17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Update the view's active selection -------------------------------
$ws.Range("AK16").Select()
